$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 527205.1
$ws.Range("I40").Value = 885
$ws.Range("J40").Value = 770122.0600000001
$ws.Range("K40").Value = 885
$ws.Range("L40").Value = 770122.0600000001
$ws.Range("M40").Value = -710
$ws.Range("N40").Value = -770472.0600000001
$ws.Range("H116").Value = 5752.8887
$ws.Range("I116").Value = 6700.1816
$ws.Range("J116").Value = 4264.2856
$ws.Range("K116").Value = 6700.1816
$ws.Range("L116").Value = 4264.2856
$ws.Range("M116").Value = -3258.1816
$ws.Range("N116").Value = -11148.2856
$ws.Range("H135").Value = 17857770
$ws.Range("I135").Value = 413.7619
$ws.Range("J135").Value = 71429840
$ws.Range("K135").Value = 3723.8571
$ws.Range("L135").Value = 642868560
$ws.Range("M135").Value = -1188.8571
$ws.Range("N135").Value = -642873630
$ws.Range("H138").Value = 4943.726
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 4943.726
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 14831.178
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -25111.178

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2613.8572
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 2613.8572
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H80").Value = 17460
$ws.Range("J80").Value = 17460
$ws.Range("L80").Value = 17460
$ws.Range("N80").Value = -19456
$ws.Range("H83").Value = 17460
$ws.Range("J83").Value = 17460
$ws.Range("L83").Value = 52380
$ws.Range("N83").Value = -62364

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2400.4
$ws.Range("I22").Value = 5250
$ws.Range("J22").Value = 500.66666
$ws.Range("K22").Value = 5250
$ws.Range("L22").Value = 500.66666
$ws.Range("M22").Value = -5077
$ws.Range("N22").Value = -846.66666
$ws.Range("H82").Value = 14730
$ws.Range("I82").Value = 7100
$ws.Range("J82").Value = 29990
$ws.Range("K82").Value = 7100
$ws.Range("L82").Value = 29990
$ws.Range("M82").Value = -6717
$ws.Range("N82").Value = -30756
$ws.Range("H85").Value = 14730
$ws.Range("I85").Value = 7100
$ws.Range("J85").Value = 29990
$ws.Range("K85").Value = 7100
$ws.Range("L85").Value = 29990
$ws.Range("M85").Value = -5774
$ws.Range("N85").Value = -32642
$ws.Range("H99").Value = 2525
$ws.Range("I99").Value = 2583.3333
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 2583.3333
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -1085.3333
$ws.Range("N99").Value = -4996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21673.477
$ws.Range("I31").Value = 23103.26
$ws.Range("J31").Value = 14654.546
$ws.Range("K31").Value = 23103.26
$ws.Range("L31").Value = 14654.546
$ws.Range("M31").Value = -22808.26
$ws.Range("N31").Value = -15244.546
$ws.Range("H34").Value = 21673.477
$ws.Range("I34").Value = 23103.26
$ws.Range("J34").Value = 14654.546
$ws.Range("K34").Value = 23103.26
$ws.Range("L34").Value = 14654.546
$ws.Range("M34").Value = -22901.26
$ws.Range("N34").Value = -15058.546
$ws.Range("H41").Value = 9225.357
$ws.Range("I41").Value = 3925
$ws.Range("J41").Value = 11345.5
$ws.Range("K41").Value = 3925
$ws.Range("L41").Value = 11345.5
$ws.Range("M41").Value = -3497
$ws.Range("N41").Value = -12201.5
$ws.Range("H50").Value = 17988
$ws.Range("J50").Value = 17985
$ws.Range("L50").Value = 17985
$ws.Range("N50").Value = -19235
$ws.Range("H51").Value = 10000
$ws.Range("I51").Value = 10000
$ws.Range("K51").Value = 10000
$ws.Range("M51").Value = -9264
$ws.Range("H58").Value = 2421.121
$ws.Range("I58").Value = 850.6111
$ws.Range("J58").Value = 4305.7334
$ws.Range("K58").Value = 850.6111
$ws.Range("L58").Value = 4305.7334
$ws.Range("M58").Value = -647.6111
$ws.Range("N58").Value = -4711.7334
$ws.Range("H59").Value = 10817.333
$ws.Range("I59").Value = 3000
$ws.Range("J59").Value = 12020
$ws.Range("K59").Value = 3000
$ws.Range("L59").Value = 12020
$ws.Range("M59").Value = -1855
$ws.Range("N59").Value = -14310
$ws.Range("H60").Value = 9394.296
$ws.Range("I60").Value = 3397.6667
$ws.Range("J60").Value = 11107.619
$ws.Range("K60").Value = 3397.6667
$ws.Range("L60").Value = 11107.619
$ws.Range("M60").Value = -2886.6667
$ws.Range("N60").Value = -12129.619
$ws.Range("H61").Value = 10000
$ws.Range("I61").Value = 10000
$ws.Range("K61").Value = 10000
$ws.Range("M61").Value = -9652
$ws.Range("H74").Value = 11214.182
$ws.Range("J74").Value = 11214.182
$ws.Range("L74").Value = 11214.182
$ws.Range("N74").Value = -12962.182
$ws.Range("H77").Value = 11214.182
$ws.Range("J77").Value = 11214.182
$ws.Range("L77").Value = 33642.546
$ws.Range("N77").Value = -42378.546
$ws.Range("H107").Value = 674.7
$ws.Range("I107").Value = 586.8333
$ws.Range("J107").Value = 806.5
$ws.Range("K107").Value = 586.8333
$ws.Range("L107").Value = 806.5
$ws.Range("M107").Value = 1333.1667
$ws.Range("N107").Value = -4646.5
$ws.Range("H136").Value = 2421.121
$ws.Range("I136").Value = 850.6111
$ws.Range("J136").Value = 4305.7334
$ws.Range("K136").Value = 2551.8333
$ws.Range("L136").Value = 12917.2002
$ws.Range("M136").Value = -1.833299999999781
$ws.Range("N136").Value = -18017.2002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 34482890
$ws.Range("I34").Value = 94
$ws.Range("J34").Value = 41666804
$ws.Range("K34").Value = 282
$ws.Range("L34").Value = 125000412
$ws.Range("M34").Value = -198
$ws.Range("N34").Value = -125000580
$ws.Range("H39").Value = 142873280
$ws.Range("I39").Value = 980
$ws.Range("J39").Value = 166685330
$ws.Range("K39").Value = 2940
$ws.Range("L39").Value = 500055990
$ws.Range("M39").Value = -2646
$ws.Range("N39").Value = -500056578
$ws.Range("H55").Value = 33337704
$ws.Range("J55").Value = 36115730
$ws.Range("L55").Value = 108347190
$ws.Range("N55").Value = -108347544
$ws.Range("H80").Value = 2347.4
$ws.Range("J80").Value = 2422.2144
$ws.Range("L80").Value = 7266.6432
$ws.Range("N80").Value = -9138.643199999999
$ws.Range("H83").Value = 2347.4
$ws.Range("J83").Value = 2422.2144
$ws.Range("L83").Value = 21799.9296
$ws.Range("N83").Value = -31159.9296
$ws.Range("H86").Value = 474.83334
$ws.Range("I86").Value = 445.7143
$ws.Range("J86").Value = 515.6
$ws.Range("K86").Value = 1337.1429
$ws.Range("L86").Value = 1546.8
$ws.Range("M86").Value = -151.1428999999998
$ws.Range("N86").Value = -3918.8
$ws.Range("H89").Value = 474.83334
$ws.Range("I89").Value = 445.7143
$ws.Range("J89").Value = 515.6
$ws.Range("K89").Value = 4011.4287
$ws.Range("L89").Value = 4640.400000000001
$ws.Range("M89").Value = 1916.5713
$ws.Range("N89").Value = -16496.4
$ws.Range("H92").Value = 820.4
$ws.Range("I92").Value = 651
$ws.Range("J92").Value = 933.3333
$ws.Range("K92").Value = 1953
$ws.Range("L92").Value = 2799.9999
$ws.Range("M92").Value = -705
$ws.Range("N92").Value = -5295.9999
$ws.Range("H103").Value = 1985.75
$ws.Range("I103").Value = 826
$ws.Range("J103").Value = 2814.1428
$ws.Range("K103").Value = 2478
$ws.Range("L103").Value = 8442.428400000001
$ws.Range("M103").Value = -1599
$ws.Range("N103").Value = -10200.4284
$ws.Range("H116").Value = 1050.2222
$ws.Range("I116").Value = 260
$ws.Range("J116").Value = 1276
$ws.Range("K116").Value = 780
$ws.Range("L116").Value = 3828
$ws.Range("M116").Value = 2662
$ws.Range("N116").Value = -10712
$ws.Range("H128").Value = 400000
$ws.Range("I128").Value = 400000
$ws.Range("K128").Value = 1200000
$ws.Range("M128").Value = -1195020

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 13259.777
$ws.Range("I43").Value = 566.6667
$ws.Range("J43").Value = 19606.334
$ws.Range("K43").Value = 566.6667
$ws.Range("L43").Value = 19606.334
$ws.Range("M43").Value = -415.6667
$ws.Range("N43").Value = -19908.334
$ws.Range("H46").Value = 5660
$ws.Range("I46").Value = 4000
$ws.Range("J46").Value = 8980
$ws.Range("K46").Value = 4000
$ws.Range("L46").Value = 8980
$ws.Range("M46").Value = -3844
$ws.Range("N46").Value = -9292
$ws.Range("H57").Value = 7966.6665
$ws.Range("J57").Value = 7966.6665
$ws.Range("L57").Value = 7966.6665
$ws.Range("N57").Value = -9606.666499999999
$ws.Range("H80").Value = 15555.556
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 25600
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 25600
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -27596
$ws.Range("H83").Value = 15555.556
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 25600
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 128000
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -137984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 970.1
$ws.Range("I22").Value = 933.5
$ws.Range("J22").Value = 1025
$ws.Range("K22").Value = 933.5
$ws.Range("L22").Value = 1025
$ws.Range("M22").Value = -638.5
$ws.Range("N22").Value = -1615
$ws.Range("H27").Value = 970.1
$ws.Range("I27").Value = 933.5
$ws.Range("J27").Value = 1025
$ws.Range("K27").Value = 933.5
$ws.Range("L27").Value = 1025
$ws.Range("M27").Value = -826.5
$ws.Range("N27").Value = -1239
$ws.Range("H46").Value = 851.8125
$ws.Range("I46").Value = 722.9
$ws.Range("J46").Value = 1066.6666
$ws.Range("K46").Value = 722.9
$ws.Range("L46").Value = 1066.6666
$ws.Range("M46").Value = -534.9
$ws.Range("N46").Value = -1442.6666

